$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row2
$ws.Range("T2").Value = 2.2
$ws.Range("U2").Value = 1.68
$ws.Range("Y2").Value = 980
$ws.Range("AD2").Value = 980
$ws.Range("AH2").Value = 980
$ws.Range("AJ2").Value = 1000
$ws.Range("AK2").Value = 1000
$ws.Range("AL2").Value = 980
$ws.Range("AN2").Value = 1000

# Row4
$ws.Range("H4").Value = 1.9
$ws.Range("L4").Value = 1.37
$ws.Range("N4").Value = 4.5
$ws.Range("Q4").Value = 1.82
$ws.Range("S4").Value = 3.1
$ws.Range("AG4").Value = 17.5

# Row5
$ws.Range("F5").Value = 3.35
$ws.Range("G5").Value = 3.4
$ws.Range("H5").Value = 2.24
$ws.Range("I5").Value = 2.26
$ws.Range("L5").Value = 1.3
$ws.Range("N5").Value = 5.6
$ws.Range("P5").Value = 2.5
$ws.Range("R5").Value = 1.61
$ws.Range("U5").Value = 2.66
$ws.Range("V5").Value = 1.79
$ws.Range("W5").Value = 1.41
$ws.Range("AA5").Value = 29
$ws.Range("AE5").Value = 19.5
$ws.Range("AG5").Value = 14.5
$ws.Range("AJ5").Value = 60
$ws.Range("AO5").Value = 12

# Row6
$ws.Range("T6").Value = 1.65
$ws.Range("X6").Value = 970
$ws.Range("Y6").Value = 970
$ws.Range("Z6").Value = 970
$ws.Range("AB6").Value = 970
$ws.Range("AE6").Value = 970
$ws.Range("AF6").Value = 970
$ws.Range("AH6").Value = 970
$ws.Range("AI6").Value = 970
$ws.Range("AJ6").Value = 970
$ws.Range("AK6").Value = 970
$ws.Range("AL6").Value = 970
$ws.Range("AN6").Value = 970
$ws.Range("AO6").Value = 970

# Row7
$ws.Range("F7").Value = 1.36
$ws.Range("J7").Value = 1.35
$ws.Range("V7").Value = 1.33

# Row8
$ws.Range("F8").Value = 1.96
$ws.Range("H8").Value = 2.16
$ws.Range("J8").Value = 3
$ws.Range("K8").Value = 970
$ws.Range("L8").Value = 1.01
$ws.Range("N8").Value = 1.33
$ws.Range("P8").Value = 1.33
$ws.Range("Q8").Value = 1.66
$ws.Range("R8").Value = 1.33
$ws.Range("S8").Value = 2.58

# Row9
$ws.Range("L9").Value = 1.01
$ws.Range("M9").Value = 1.12
$ws.Range("N9").Value = 2.32
$ws.Range("O9").Value = 1.56
$ws.Range("P9").Value = 1.46
$ws.Range("Q9").Value = 2.94
$ws.Range("R9").Value = 1.15
$ws.Range("S9").Value = 5.5
$ws.Range("T9").Value = 1.89
$ws.Range("U9").Value = 1.47
$ws.Range("V9").Value = 1.3
$ws.Range("W9").Value = 1.67
$ws.Range("X9").Value = 9.8
$ws.Range("Y9").Value = 14.5
$ws.Range("Z9").Value = 970
$ws.Range("AA9").Value = 1000
$ws.Range("AB9").Value = 9.6
$ws.Range("AC9").Value = 10
$ws.Range("AD9").Value = 28
$ws.Range("AE9").Value = 100
$ws.Range("AF9").Value = 19
$ws.Range("AG9").Value = 18.5
$ws.Range("AH9").Value = 38
$ws.Range("AI9").Value = 1000
$ws.Range("AJ9").Value = 50
$ws.Range("AK9").Value = 55
$ws.Range("AL9").Value = 100
$ws.Range("AM9").Value = 1000
$ws.Range("AN9").Value = 1000
$ws.Range("AO9").Value = 1000

# Row10
$ws.Range("F10").Value = 2.68
$ws.Range("G10").Value = 3.8
$ws.Range("H10").Value = 2.48
$ws.Range("I10").Value = 3.3
$ws.Range("J10").Value = 2.66
$ws.Range("K10").Value = 3.6
$ws.Range("L10").Value = 1.01
$ws.Range("M10").Value = 1.01
$ws.Range("N10").Value = 1.4
$ws.Range("O10").Value = 1.01
$ws.Range("R10").Value = 1.14
$ws.Range("S10").Value = 4.5
$ws.Range("T10").Value = 1.78
$ws.Range("U10").Value = 1.55
$ws.Range("V10").Value = 1.43
$ws.Range("W10").Value = 1.41
$ws.Range("X10").Value = 11.5
$ws.Range("Y10").Value = 970
$ws.Range("Z10").Value = 970
$ws.Range("AA10").Value = 1000
$ws.Range("AB10").Value = 970
$ws.Range("AC10").Value = 970
$ws.Range("AD10").Value = 970
$ws.Range("AE10").Value = 1000
$ws.Range("AF10").Value = 970
$ws.Range("AG10").Value = 970
$ws.Range("AH10").Value = 1000
$ws.Range("AI10").Value = 1000
$ws.Range("AJ10").Value = 1000
$ws.Range("AK10").Value = 1000
$ws.Range("AL10").Value = 1000
$ws.Range("AM10").Value = 1000
$ws.Range("AN10").Value = 1000
$ws.Range("AO10").Value = 1000

# Row11
$ws.Range("L11").Value = 1.01
$ws.Range("M11").Value = 1.01
$ws.Range("N11").Value = 1.02
$ws.Range("O11").Value = 1.42
$ws.Range("R11").Value = 1.08
$ws.Range("S11").Value = 1.01
$ws.Range("T11").Value = 1.01
$ws.Range("U11").Value = 1.01
$ws.Range("V11").Value = 1.01
$ws.Range("W11").Value = 1.01
$ws.Range("X11").Value = 1000
$ws.Range("Y11").Value = 1000
$ws.Range("Z11").Value = 1000
$ws.Range("AA11").Value = 1000
$ws.Range("AB11").Value = 1000
$ws.Range("AC11").Value = 1000
$ws.Range("AD11").Value = 1000
$ws.Range("AE11").Value = 1000
$ws.Range("AF11").Value = 1000
$ws.Range("AG11").Value = 1000
$ws.Range("AH11").Value = 1000
$ws.Range("AI11").Value = 1000
$ws.Range("AJ11").Value = 1000
$ws.Range("AK11").Value = 1000
$ws.Range("AL11").Value = 1000
$ws.Range("AM11").Value = 1000
$ws.Range("AN11").Value = 1000
$ws.Range("AO11").Value = 1000

# Row12
$ws.Range("H12").Value = 4.9
$ws.Range("J12").Value = 2.86
$ws.Range("K12").Value = 3.25

# Row14
$ws.Range("H14").Value = 4.5
$ws.Range("J14").Value = 3.1
